$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
$ws.Range("A1").Value = "test"
